$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Melissa's email (C2): melissa@gmail.com -> melissa@renting.com ---
# Remove the old hyperlink first so we don't leave a stale mailto: relationship,
# then re-create it pointing at the new address.
$ws.Range("C2").Hyperlinks.Delete()
$ws.Range("C2").Value = "melissa@renting.com"
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:melissa@renting.com", "", "", "melissa@renting.com")
$ws.Range("C2").Style = "Hyperlink"

# --- Add a new agent row (row 3): Jack ---
$ws.Range("A3").Value = 789101
$ws.Range("B3").Value = "Jack"
$ws.Range("C3").Value = "jack@renting.com"
$ws.Range("D3").Value = 2158889999
$ws.Range("E3").Value = "hbye123"
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:jack@renting.com", "", "", "jack@renting.com")
$ws.Range("C3").Style = "Hyperlink"

# --- Resize email column to fit the new, longer addresses (best-fit width) ---
$ws.Columns("C").ColumnWidth = 18

# --- Move the active selection down to the next empty row, like a user about
# to enter a new record ---
$ws.Range("A4").Select()

Write-Output "done"
